$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string: force text format
# (matching original style where all data cells are stored as text)
# so Excel does not auto-convert them into real numbers.
$textFormatCells = @(
    "D6",
    "D9",
    "D13",
    "D17",
    "D19",
    "D22",
    "D23",
    "D24",
    "D30",
    "D31",
    "D33",
    "D36",
    "D43",
    "D44",
    "D46",
    "D48"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Range("D2").Value = '60.771.79'
$ws.Range("D3").Value = '2.907.64'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("D6").Value = '144.25'
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").Value = '6.87'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("E11").Value = '  -2.32%  '
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '33.34'
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").Value = '3.390.39'
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '60.708.94'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '6.66'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '2.906.24'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").Value = '430.96'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("E20").Value = '  -2.12%  '
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").Value = '7.07'
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").Value = '81.24'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").Value = '10.84'
$ws.Range("E24").Value = '  +0.80%  '
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  +4.36%  '
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").Value = '6.94'
$ws.Range("E30").Value = '  -3.75%  '
$ws.Range("D31").Value = '26.48'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '0.0₃0855'
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("D36").Value = '5.60'
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("E39").Value = '  -5.03%  '
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("E42").Value = '  -5.81%  '
$ws.Range("D43").Value = '374.56'
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0343'
$ws.Range("E44").Value = '  -2.76%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.693.36'
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("D46").Value = '133.72'
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '23.70'
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("E49").Value = '  -0.99%  '
$ws.Range("E50").Value = '  -3.27%  '
$ws.Range("E51").Value = '  -1.05%  '

Write-Output "Applied cryptos list update"
